# Auto-generated Excel COM-interop script to apply cell value updates
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 8932.666999999999
$ws.Range("I86").Value = 9150
$ws.Range("J86").Value = 8498
$ws.Range("K86").Value = 9150
$ws.Range("L86").Value = 8498
$ws.Range("M86").Value = -8027
$ws.Range("N86").Value = -10744
$ws.Range("H89").Value = 8932.666999999999
$ws.Range("I89").Value = 9150
$ws.Range("J89").Value = 8498
$ws.Range("K89").Value = 45750
$ws.Range("L89").Value = 42490
$ws.Range("M89").Value = -40134
$ws.Range("N89").Value = -53722
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").Value = $null
$ws.Range("H112").Value = 2464.8462
$ws.Range("J112").Value = 2999.875
$ws.Range("L112").Value = 8999.625
$ws.Range("N112").Value = -11215.625
$ws.Range("H113").Value = 9249.5
$ws.Range("I113").Value = 9249.5
$ws.Range("K113").Value = 9249.5
$ws.Range("M113").Value = -5995.5
$ws.Range("H132").Value = 2478.923
$ws.Range("I132").Value = 1182.7
$ws.Range("K132").Value = 3548.1
$ws.Range("M132").Value = -1018.1

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8355.6
$ws.Range("I32").Value = 8355.6
$ws.Range("K32").Value = 8355.6
$ws.Range("M32").Value = -8068.6
$ws.Range("H45").Value = 24599.334
$ws.Range("I45").Value = 7798.8
$ws.Range("J45").Value = 45600
$ws.Range("K45").Value = 7798.8
$ws.Range("L45").Value = 45600
$ws.Range("M45").Value = -7421.8
$ws.Range("N45").Value = -46354
$ws.Range("H74").Value = 6205.7
$ws.Range("I74").Value = 6795.222
$ws.Range("J74").Value = 900
$ws.Range("K74").Value = 6795.222
$ws.Range("L74").Value = 900
$ws.Range("M74").Value = -5921.222
$ws.Range("N74").Value = -2648
$ws.Range("H77").Value = 6205.7
$ws.Range("I77").Value = 6795.222
$ws.Range("J77").Value = 900
$ws.Range("K77").Value = 33976.11
$ws.Range("L77").Value = 4500
$ws.Range("M77").Value = -29608.11
$ws.Range("N77").Value = -13236
$ws.Range("H110").Value = 1529.0769
$ws.Range("I110").Value = 1529.0769
$ws.Range("K110").Value = 1529.0769
$ws.Range("M110").Value = 515.9231
$ws.Range("H132").Value = 4244.2856
$ws.Range("I132").Value = 3904
$ws.Range("K132").Value = 11712
$ws.Range("M132").Value = -9182

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 562.5
$ws.Range("I7").Value = 1000
$ws.Range("K7").Value = 1000
$ws.Range("M7").Value = -887
$ws.Range("H107").Value = 1403.6666
$ws.Range("I107").Value = 1155.5
$ws.Range("K107").Value = 1155.5
$ws.Range("M107").Value = 764.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3401
$ws.Range("I31").Value = 3370.3333
$ws.Range("J31").Value = 3416.3333
$ws.Range("K31").Value = 3370.3333
$ws.Range("L31").Value = 3416.3333
$ws.Range("M31").Value = -3075.3333
$ws.Range("N31").Value = -4006.3333
$ws.Range("H34").Value = 3401
$ws.Range("I34").Value = 3370.3333
$ws.Range("J34").Value = 3416.3333
$ws.Range("K34").Value = 3370.3333
$ws.Range("L34").Value = 3416.3333
$ws.Range("M34").Value = -3168.3333
$ws.Range("N34").Value = -3820.3333
$ws.Range("H58").Value = 3337.0833
$ws.Range("I58").Value = 3597
$ws.Range("K58").Value = 3597
$ws.Range("M58").Value = -3394
$ws.Range("H62").Value = 2124.5
$ws.Range("I62").Value = 1833
$ws.Range("J62").Value = 2999
$ws.Range("K62").Value = 1833
$ws.Range("L62").Value = 2999
$ws.Range("M62").Value = -1209
$ws.Range("N62").Value = -4247
$ws.Range("H65").Value = 2124.5
$ws.Range("I65").Value = 1833
$ws.Range("J65").Value = 2999
$ws.Range("K65").Value = 9165
$ws.Range("L65").Value = 14995
$ws.Range("M65").Value = -6045
$ws.Range("N65").Value = -21235
$ws.Range("H107").Value = 101522.8
$ws.Range("I107").Value = 168371.5
$ws.Range("K107").Value = 168371.5
$ws.Range("M107").Value = -166451.5
$ws.Range("H134").Value = 6371
$ws.Range("I134").Value = 9150
$ws.Range("K134").Value = 27450
$ws.Range("M134").Value = -24915
$ws.Range("H136").Value = 3337.0833
$ws.Range("I136").Value = 3597
$ws.Range("K136").Value = 10791
$ws.Range("M136").Value = -8241

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1659.091
$ws.Range("I131").Value = 1562.5
$ws.Range("J131").Value = 1714.2858
$ws.Range("K131").Value = 4687.5
$ws.Range("L131").Value = 5142.857400000001
$ws.Range("M131").Value = 352.5
$ws.Range("N131").Value = -15222.8574

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 358.16666
$ws.Range("J2").Value = 566.3333
$ws.Range("L2").Value = 566.3333
$ws.Range("N2").Value = -792.3333
$ws.Range("H122").Value = 2271.1177
$ws.Range("I122").Value = 2369.923
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 7109.768999999999
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -4659.768999999999
$ws.Range("N122").Value = -10750
$ws.Range("H126").Value = 3994
$ws.Range("H132").Value = 3570.3635
$ws.Range("I132").Value = 2659.5
$ws.Range("J132").Value = 5999.3335
$ws.Range("K132").Value = 7978.5
$ws.Range("L132").Value = 17998.0005
$ws.Range("M132").Value = -5448.5
$ws.Range("N132").Value = -23058.0005

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8228.611000000001
$ws.Range("I22").Value = 7928.727
$ws.Range("K22").Value = 7928.727
$ws.Range("M22").Value = -7633.727
$ws.Range("H27").Value = 8228.611000000001
$ws.Range("I27").Value = 7928.727
$ws.Range("K27").Value = 7928.727
$ws.Range("M27").Value = -7821.727
$ws.Range("H132").Value = 3095.8
$ws.Range("I132").Value = 1851.4286
$ws.Range("K132").Value = 5554.2858
$ws.Range("M132").Value = -3024.2858

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3628.5715
$ws.Range("I126").Value = 3628.5715
$ws.Range("K126").Value = 10885.7145
$ws.Range("M126").Value = -8415.7145
$ws.Range("H132").Value = 1665.4445
$ws.Range("I132").Value = 921.53845
$ws.Range("K132").Value = 2764.61535
$ws.Range("M132").Value = -234.61535
